# Update "想去人数" (interest count) figures across the workbook sheets.
# Sheet 1: 展览 (Exhibitions)
# Sheet 2: 演出 (Performances)
# Sheet 3: 本地生活 (Local Life)
# Sheet 4: 全部类型 (All Types) - combined listing

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet 1) ---
$ws1.Range("F4").Value  = 1178
$ws1.Range("F5").Value  = 219
$ws1.Range("F6").Value  = 2794
$ws1.Range("F8").Value  = 711
$ws1.Range("F9").Value  = 126
$ws1.Range("F11").Value = 205
$ws1.Range("F12").Value = 713
$ws1.Range("F13").Value = 118
$ws1.Range("F14").Value = 138
$ws1.Range("F15").Value = 1856
$ws1.Range("F18").Value = 209
$ws1.Range("F19").Value = 261

# --- 演出 (sheet 2) ---
$ws2.Range("F6").Value  = 20
$ws2.Range("F7").Value  = 27
$ws2.Range("F10").Value = 62
$ws2.Range("F12").Value = 59
$ws2.Range("F22").Value = 47
$ws2.Range("F23").Value = 37

# --- 本地生活 (sheet 3) ---
$ws3.Range("F2").Value = 6371
$ws3.Range("F3").Value = 800
$ws3.Range("F4").Value = 2038
$ws3.Range("F5").Value = 276

# --- 全部类型 (sheet 4) ---
$ws4.Range("F2").Value  = 6371
$ws4.Range("F3").Value  = 800
$ws4.Range("F4").Value  = 2038
$ws4.Range("F5").Value  = 276
$ws4.Range("F12").Value = 1178
$ws4.Range("F13").Value = 219
$ws4.Range("F14").Value = 20
$ws4.Range("F15").Value = 27
$ws4.Range("F17").Value = 2794
$ws4.Range("F20").Value = 62
$ws4.Range("F22").Value = 59
$ws4.Range("F23").Value = 711
$ws4.Range("F24").Value = 126
$ws4.Range("F27").Value = 205
$ws4.Range("F28").Value = 713
$ws4.Range("F29").Value = 118
$ws4.Range("F30").Value = 138
$ws4.Range("F32").Value = 1856
$ws4.Range("F37").Value = 209
$ws4.Range("F43").Value = 47
$ws4.Range("F44").Value = 37
$ws4.Range("F45").Value = 261
